$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.712771666666666
$ws.Range("H2").Value = 17.138315
$ws.Range("I2").Value = 0.1683613830606884
$ws.Range("J2").Value = 0.1683613830606885
$ws.Range("M2").Value = 0.8272046666666667
$ws.Range("N2").Value = 2.481614
$ws.Range("O2").Value = 0.2219283187021488
$ws.Range("P2").Value = 0.2219283187021487
$ws.Range("Q2").Value = 4.725631382267777
$ws.Range("R2").Value = 42.53068244041
$ws.Range("S2").Value = 0.03736415867702701
$ws.Range("T2").Value = 0.03736415867702701
$ws.Range("G3").Value = 5.712771666666666
$ws.Range("H3").Value = 17.138315
$ws.Range("I3").Value = 0.1683613830606884
$ws.Range("J3").Value = 0.1683613830606885
$ws.Range("O3").Value = 0.1711024058911549
$ws.Range("P3").Value = 0.1711024058911549
$ws.Range("Q3").Value = 3.643369641104444
$ws.Range("R3").Value = 32.79032676993999
$ws.Range("S3").Value = 0.02880703770084611
$ws.Range("T3").Value = 0.02880703770084612
$ws.Range("G4").Value = 5.712771666666666
$ws.Range("H4").Value = 17.138315
$ws.Range("I4").Value = 0.1683613830606884
$ws.Range("J4").Value = 0.1683613830606885
$ws.Range("M4").Value = 2.262387333333333
$ws.Range("N4").Value = 6.787162
$ws.Range("O4").Value = 0.6069692754066963
$ws.Range("P4").Value = 0.6069692754066964
$ws.Range("Q4").Value = 12.92450225689222
$ws.Range("R4").Value = 116.32052031203
$ws.Range("S4").Value = 0.1021901866828153
$ws.Range("T4").Value = 0.1021901866828153
$ws.Range("I5").Value = 0.4370667227533506
$ws.Range("J5").Value = 0.4370667227533506
$ws.Range("M5").Value = 0.8272046666666667
$ws.Range("N5").Value = 2.481614
$ws.Range("O5").Value = 0.2219283187021488
$ws.Range("P5").Value = 0.2219283187021487
$ws.Range("Q5").Value = 12.26775513268178
$ws.Range("R5").Value = 110.409796194136
$ws.Range("S5").Value = 0.09699748294130928
$ws.Range("T5").Value = 0.09699748294130926
$ws.Range("I6").Value = 0.4370667227533506
$ws.Range("J6").Value = 0.4370667227533506
$ws.Range("O6").Value = 0.1711024058911549
$ws.Range("P6").Value = 0.1711024058911549
$ws.Range("R6").Value = 85.123799762224
$ws.Range("S6").Value = 0.07478316779806064
$ws.Range("T6").Value = 0.07478316779806064
$ws.Range("I7").Value = 0.4370667227533506
$ws.Range("J7").Value = 0.4370667227533506
$ws.Range("M7").Value = 2.262387333333333
$ws.Range("N7").Value = 6.787162
$ws.Range("O7").Value = 0.6069692754066963
$ws.Range("P7").Value = 0.6069692754066964
$ws.Range("Q7").Value = 33.55205179445422
$ws.Range("R7").Value = 301.968466150088
$ws.Range("S7").Value = 0.2652860720139806
$ws.Range("T7").Value = 0.2652860720139807
$ws.Range("G8").Value = 13.388457
$ws.Range("H8").Value = 40.165371
$ws.Range("I8").Value = 0.3945718941859609
$ws.Range("J8").Value = 0.3945718941859609
$ws.Range("M8").Value = 0.8272046666666667
$ws.Range("N8").Value = 2.481614
$ws.Range("O8").Value = 0.2219283187021488
$ws.Range("P8").Value = 0.2219283187021487
$ws.Range("Q8").Value = 11.074994109866
$ws.Range("R8").Value = 99.674946988794
$ws.Range("S8").Value = 0.08756667708381247
$ws.Range("T8").Value = 0.08756667708381245
$ws.Range("G9").Value = 13.388457
$ws.Range("H9").Value = 40.165371
$ws.Range("I9").Value = 0.3945718941859609
$ws.Range("J9").Value = 0.3945718941859609
$ws.Range("O9").Value = 0.1711024058911549
$ws.Range("P9").Value = 0.1711024058911549
$ws.Range("Q9").Value = 8.538604485044
$ws.Range("R9").Value = 76.847440365396
$ws.Range("S9").Value = 0.06751220039224809
$ws.Range("T9").Value = 0.06751220039224809
$ws.Range("G10").Value = 13.388457
$ws.Range("H10").Value = 40.165371
$ws.Range("I10").Value = 0.3945718941859609
$ws.Range("J10").Value = 0.3945718941859609
$ws.Range("M10").Value = 2.262387333333333
$ws.Range("N10").Value = 6.787162
$ws.Range("O10").Value = 0.6069692754066963
$ws.Range("P10").Value = 0.6069692754066964
$ws.Range("Q10").Value = 30.289875529678
$ws.Range("R10").Value = 272.608879767102
$ws.Range("S10").Value = 0.2394930167099004
$ws.Range("T10").Value = 0.2394930167099004
Write-Output "Updated TPM-derived values in G2:T10"
